$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9235045313835144
$ws.Range("B1").Value = 0.7597807645797729
$ws.Range("C1").Value = 0.5873104333877563
$ws.Range("D1").Value = 0.5677320957183838
$ws.Range("E1").Value = 0.6123800277709961
